# "Generate Report for Handback" — mark a.md/b.md as handed back (in sync
# with en-US) for both locales, recording the handback target file/time.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: both locale columns flip from "Ready for handoff" to
# "Handed back: in sync with en-US" for the a.md / b.md rows.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusHandedBack
$overview.Range("C2").Value = $statusHandedBack
$overview.Range("B3").Value = $statusHandedBack
$overview.Range("C3").Value = $statusHandedBack

# ---------------------------------------------------------------------
# Per-locale detail sheets: status flips the same way, and a new
# "Latest Target File" / "Latest Handback File" pair + "Latest Handback
# DateTime" get filled in for the two content rows (a.md, b.md).
# ---------------------------------------------------------------------
function Set-HandbackRow($ws, $row, $targetFileName, $targetFileUrl, $handbackDateTime) {
    $ws.Cells.Item($row, 2).Value = $statusHandedBack

    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Value = "a.md"
    $ws.Hyperlinks.Add($eCell, "https://github.com/OpenLocalizationTest/oltest/blob/7fabfbbcd3d949f55cc06aa56ace591e0f4cefbc/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null

    $fCell = $ws.Cells.Item($row, 6)
    $fCell.Value = $targetFileName
    $ws.Hyperlinks.Add($fCell, $targetFileUrl, [Type]::Missing, [Type]::Missing, $targetFileName) | Out-Null

    $ws.Cells.Item($row, 7).Value = $handbackDateTime
}

# zh-cn sheet
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ef69fd6fc186252e3edf41957032c35d91a39c27/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
Set-HandbackRow $zhcn 2 $zhXlfName $zhXlfUrl "2016-02-25 05:43:44"
Set-HandbackRow $zhcn 3 $zhXlfName $zhXlfUrl "2016-02-25 05:43:44"

# de-de sheet
$dede = $wb.Worksheets.Item("de-de")
$deXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0f82380ca4102338ac49fa54b069b1a2979c6fc9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
Set-HandbackRow $dede 2 $deXlfName $deXlfUrl "2016-02-25 05:44:04"
Set-HandbackRow $dede 3 $deXlfName $deXlfUrl "2016-02-25 05:44:04"
